# Update cryptos list values (Price column D, Volume(1h) column E)
# D-column values are plain numeric-looking strings in the source data; force
# Text number format before assignment so COM does not silently coerce them to
# doubles (losing exact text / precision), then clear the format again so the
# cell keeps its original (default/no) style, matching the source workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '43.072.16'
$r.ClearFormats()
$ws.Range("E2").Value = '  +0.16%  '
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '2.300.73'
$r.ClearFormats()
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  -0.05%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '300.01'
$r.ClearFormats()
$ws.Range("E5").Value = '  -0.26%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '97.70'
$r.ClearFormats()
$ws.Range("E6").Value = '  -1.69%  '
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.521'
$r.ClearFormats()
$ws.Range("E7").Value = '  +3.51%  '
$ws.Range("E8").Value = '  -0.03%  '
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.515'
$r.ClearFormats()
$ws.Range("E9").Value = '  +0.84%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '36.21'
$r.ClearFormats()
$ws.Range("E10").Value = '  +0.10%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.0793'
$r.ClearFormats()
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("E12").Value = '  +0.51%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '17.77'
$r.ClearFormats()
$ws.Range("E13").Value = '  -2.75%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '6.88'
$r.ClearFormats()
$ws.Range("E14").Value = '  -0.95%  '
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '2.655.88'
$r.ClearFormats()
$ws.Range("E15").Value = '  +0.05%  '
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '2.303.71'
$r.ClearFormats()
$ws.Range("E16").Value = '  +0.91%  '
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '0.788'
$r.ClearFormats()
$ws.Range("E17").Value = '  -1.39%  '
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '42.934.54'
$r.ClearFormats()
$ws.Range("E18").Value = '  +0.07%  '
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '13.01'
$r.ClearFormats()
$ws.Range("E19").Value = '  +3.83%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '0.0₃0911'
$r.ClearFormats()
$ws.Range("E20").Value = '  +0.87%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '6.13'
$r.ClearFormats()
$ws.Range("E21").Value = '  +0.29%  '
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '68.21'
$r.ClearFormats()
$ws.Range("E22").Value = '  +0.69%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '237.82'
$r.ClearFormats()
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("E26").Value = '  -0.49%  '
$ws.Range("E27").Value = '  -0.21%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '24.99'
$r.ClearFormats()
$ws.Range("E28").Value = '  -0.01%  '
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '2.05'
$r.ClearFormats()
$ws.Range("E30").Value = '  +0.17%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '163.20'
$r.ClearFormats()
$ws.Range("E31").Value = '  -2.44%  '
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '33.01'
$r.ClearFormats()
$ws.Range("E32").Value = '  -4.15%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  +2.24%  '
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '18.14'
$r.ClearFormats()
$ws.Range("E35").Value = '  +2.82%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '4.74'
$r.ClearFormats()
$ws.Range("E36").Value = '  +2.11%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '2.41'
$r.ClearFormats()
$ws.Range("E37").Value = '  +0.25%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.0697'
$r.ClearFormats()
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("E41").Value = '  +1.48%  '
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '2.76'
$r.ClearFormats()
$ws.Range("E42").Value = '  -1.90%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '2.009.50'
$r.ClearFormats()
$ws.Range("E43").Value = '  +2.06%  '
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("E45").Value = '  -3.89%  '
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '10.26'
$r.ClearFormats()
$ws.Range("E46").Value = '  +0.88%  '
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '17.41'
$r.ClearFormats()
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("E48").Value = '  -1.65%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '54.28'
$r.ClearFormats()
$ws.Range("E49").Value = '  -2.04%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '2.528.65'
$r.ClearFormats()
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("E51").Value = '  -0.45%  '
